# Daily attendance processing - 2026-01-19 19:14:50
# Reorders the "Recorded By" text in column G from "System, <email>" to
# "<email>, System" for the rows that were touched by today's processing
# run. The cells touched form contiguous runs (by row) of the exact text
# "System, dnasr281@gmail.com"; within each such run, every row is
# rewritten except the final row of runs that are 3 rows or longer
# (matching the set of rows the daily job actually re-saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$OLD_TEXT = "System, dnasr281@gmail.com"
$NEW_TEXT = "dnasr281@gmail.com, System"

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count
$col = 7  # Column G - "Recorded By"

$runStart = -1
$runLength = 0

for ($r = 1; $r -le ($lastRow + 1); $r++) {
    $isMatch = $false
    if ($r -le $lastRow) {
        $text = $ws.Cells.Item($r, $col).Text
        if ($text -eq $OLD_TEXT) {
            $isMatch = $true
        }
    }

    if ($isMatch) {
        if ($runStart -eq -1) {
            $runStart = $r
        }
        $runLength = $runLength + 1
    } else {
        if ($runLength -gt 0) {
            $runEnd = $runStart + $runLength - 1
            if ($runLength -eq 2) {
                $lastToChange = $runEnd
            } elseif ($runLength -ge 3) {
                $lastToChange = $runEnd - 1
            } else {
                $lastToChange = $runStart - 1
            }

            for ($rr = $runStart; $rr -le $lastToChange; $rr++) {
                $ws.Cells.Item($rr, $col).Value = $NEW_TEXT
            }
        }
        $runStart = -1
        $runLength = 0
    }
}
